# Applies the Tue Feb 27 06:28:56 UTC 2024 "Updated cryptos list" GitHub Actions refresh
# to the cryptos worksheet: new Price/Volume(1h) figures for most rows, plus a few rows
# whose Coin/Link/Price/Volume were re-ranked (rows 16-17, 48-50).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '56.262.27'
$ws.Range("E2").Value = '  +9.20%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '3.228.76'
$ws.Range("E3").Value = '  +4.09%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  -0.05%  '

# Row 5: BNB
$ws.Range("D5").Value = '''399.43'
$ws.Range("D5").Style = "Normal"

# Row 6: Solana
$ws.Range("D6").Value = '''110.71'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.43%  '

# Row 7: XRP
$ws.Range("E7").Value = '  +2.67%  '

# Row 8: USDC
$ws.Range("E8").Value = '  -0.05%  '

# Row 9: Cardano
$ws.Range("D9").Value = '''0.626'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.84%  '

# Row 10: Avalanche
$ws.Range("D10").Value = '''39.53'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.38%  '

# Row 11: Dogecoin
$ws.Range("D11").Value = '''0.0909'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.10%  '

# Row 12: TRON
$ws.Range("E12").Value = '  +2.16%  '

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = '3.736.09'
$ws.Range("E13").Value = '  +3.88%  '

# Row 14: Polkadot
$ws.Range("D14").Value = '''8.16'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.31%  '

# Row 15: Chainlink
$ws.Range("D15").Value = '''19.10'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.02%  '

# Row 16: WrappedEther
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.226.82'
$ws.Range("E16").Value = '  +4.30%  '

# Row 17: Polygon
$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D17").Value = '''1.06'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.75%  '

# Row 18: Uniswap
$ws.Range("D18").Value = '''10.72'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.75%  '

# Row 19: WrappedBTC
$ws.Range("D19").Value = '56.064.88'
$ws.Range("E19").Value = '  +8.67%  '

# Row 20: ImmutableX
$ws.Range("E20").Value = '  +2.61%  '

# Row 21: ShibaInu
$ws.Range("E21").Value = '  +6.60%  '

# Row 22: InternetComputer(DFINITY)
$ws.Range("E22").Value = '  +4.86%  '

# Row 23: BitcoinCash
$ws.Range("D23").Value = '''303.80'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +13.89%  '

# Row 24: Litecoin
$ws.Range("D24").Value = '''75.47'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +7.75%  '

# Row 25: PancakeSwap
$ws.Range("D25").Value = '''3.25'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.09%  '

# Row 26: Filecoin
$ws.Range("E26").Value = '  +1.65%  '

# Row 27: EthereumClassic
$ws.Range("D27").Value = '''28.30'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.61%  '

# Row 28: RenderToken
$ws.Range("D28").Value = '''7.49'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.87%  '

# Row 29: Kaspa
$ws.Range("E29").Value = '  +4.52%  '

# Row 30: Dai
$ws.Range("E30").Value = '  -0.07%  '

# Row 31: Hedera
$ws.Range("E31").Value = '  +4.93%  '

# Row 32: Cosmos
$ws.Range("D32").Value = '''11.24'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +7.70%  '

# Row 33: VeChain
$ws.Range("D33").Value = '''0.0493'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.61%  '

# Row 34: InjectiveProtocol
$ws.Range("D34").Value = '''36.28'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.53%  '

# Row 35: Toncoin
$ws.Range("D35").Value = '''2.11'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.25%  '

# Row 36: OKB
$ws.Range("D36").Value = '''51.43'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.82%  '

# Row 37: Stacks
$ws.Range("D37").Value = '''3.12'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +24.27%  '

# Row 38: FirstDigitalUSD
$ws.Range("E38").Value = '  -0.04%  '

# Row 39: LidoDAOToken
$ws.Range("E39").Value = '  +4.13%  '

# Row 40: Monero
$ws.Range("D40").Value = '''135.89'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.39%  '

# Row 41: ARBITRUM
$ws.Range("E41").Value = '  +3.18%  '

# Row 42: NEARProtocol
$ws.Range("E42").Value = '  +6.12%  '

# Row 43: Celestia
$ws.Range("E43").Value = '  +4.24%  '

# Row 44: Stellar
$ws.Range("E44").Value = '  +3.54%  '

# Row 45: TheGraph
$ws.Range("D45").Value = '''0.286'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.20%  '

# Row 46: EnergySwap
$ws.Range("D46").Value = '''22.31'
$ws.Range("D46").Style = "Normal"

# Row 47: ThetaToken
$ws.Range("E47").Value = '  +47.74%  '

# Row 48: WEMIXToken
$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").Value = '''2.12'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.70%  '

# Row 49: ApeXProtocol
$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").Value = '''2.49'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.16%  '

# Row 50: Maker
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '2.144.42'
$ws.Range("E50").Value = '  +3.36%  '

# Row 51: BEAM
$ws.Range("D51").Value = '''0.0363'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +11.06%  '
